# Regenerate the handoff report: update the GUID-based file names,
# xliff file names, and handoff/generate timestamps across the
# Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "4f86372d-57ef-4c39-9246-ea8317c06185"
$newGuid = "9d613ea1-a056-42a2-82ce-6bf0c780f2a4"

$oldZhXlf = "$oldGuid.149b321bbe1d45dc650932e9fae1d24fcaf90774.zh-cn.xlf"
$newZhXlf = "$newGuid.812b40031a2078c6f0852f659780606da3ad265e.zh-cn.xlf"

$oldDeXlf = "$oldGuid.149b321bbe1d45dc650932e9fae1d24fcaf90774.de-de.xlf"
$newDeXlf = "$newGuid.812b40031a2078c6f0852f659780606da3ad265e.de-de.xlf"

# Note: the hyperlink objects themselves (Overview!B2, zh-cn!A2, de-de!A2)
# keep pointing at the same target URL/relationship (r:id is unchanged in
# the diff) - only their cached "display" text tracks the cell text, which
# Excel keeps in sync with the cell's own value automatically. So we only
# need to update the cell values below; we deliberately avoid touching the
# Hyperlinks collection directly (Add/Delete) to avoid creating a second,
# redundant hyperlink entry alongside the original one.

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 10:38:23"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-22 10:38:18"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
